{"js": "const pairs = [\n  [\"2023-05-01 Monday\", \"2023-05-02 Tuesday\"],\n  [\"84+10=94\", \"50+0=50\"],\n  [\"56-36=20\", \"4+60=64\"],\n  [\"9+5=14\", \"36-17=19\"],\n  [\"38+3=41\", \"77+18=95\"],\n  [\"85-0=85\", \"33+0=33\"],\n  [\"50-39=11\", \"70+18=88\"],\n  [\"0+58=58\", \"85-59=26\"],\n  [\"68-11=57\", \"59-29=30\"],\n  [\"23+47=70\", \"72-28=44\"],\n  [\"44+11=55\", \"10+62=72\"],\n  [\"79-0=79\", \"83-9=74\"],\n  [\"25+73=98\", \"39+50=89\"],\n  [\"46-31=15\", \"67+14=81\"],\n  [\"27+46=73\", \"76+22=98\"],\n  [\"73-52=21\", \"47-18=29\"],\n  [\"94-60=34\", \"59-45=14\"],\n  [\"40-0=40\", \"80-14=66\"],\n  [\"2+69=71\", \"81-15=66\"],\n  [\"15+56=71\", \"58-14=44\"],\n  [\"8+56=64\", \"65-16=49\"],\n  [\"58-57=1\", \"86-19=67\"],\n  [\"9-6=3\", \"40+43=83\"],\n  [\"10+43=53\", \"68-24=44\"],\n  [\"36+25=61\", \"54+13=67\"],\n  [\"85+9=94\", \"28+40=68\"],\n  [\"71-65=6\", \"1+88=89\"],\n  [\"86-35=51\", \"0+16=16\"],\n  [\"94-57=37\", \"44-35=9\"],\n  [\"83-26=57\", \"94+0=94\"],\n  [\"70-15=55\", \"3+86=89\"],\n  [\"70-38=32\", \"49-32=17\"],\n  [\"50+14=64\", \"64-51=13\"],\n  [\"79-15=64\", \"25+44=69\"],\n  [\"97-37=60\", \"71-54=17\"],\n  [\"42+50=92\", \"16+71=87\"],\n  [\"36+8=44\", \"30+25=55\"],\n  [\"39+24=63\", \"36-33=3\"],\n  [\"13+34=47\", \"98-55=43\"],\n  [\"31+12=43\", \"64-43=21\"],\n  [\"44-9=35\", \"74-11=63\"],\n  [\"7+53=60\", \"43-12=31\"],\n  [\"68-31=37\", \"47+46=93\"],\n  [\"72+19=91\", \"86-2=84\"],\n  [\"39+31=70\", \"64-63=1\"],\n  [\"53-34=19\", \"43+34=77\"],\n  [\"34+43=77\", \"55-12=43\"],\n  [\"83+3=86\", \"33-23=10\"],\n  [\"86-70=16\", \"19-19=0\"],\n  [\"36-3=33\", \"33+8=41\"],\n  [\"45-32=13\", \"74-3=71\"],\n  [\"55-50=5\", \"10+21=31\"],\n  [\"11+57=68\", \"59-49=10\"],\n  [\"49-44=5\", \"58-51=7\"],\n  [\"45-28=17\", \"87-86=1\"],\n  [\"82-52=30\", \"54-49=5\"],\n  [\"55-46=9\", \"28+47=75\"],\n  [\"61-29=32\", \"96-18=78\"],\n  [\"11+5=16\", \"43+18=61\"],\n  [\"26-24=2\", \"85-19=66\"],\n  [\"44-27=17\", \"66-19=47\"],\n  [\"12+62=74\", \"81-45=36\"],\n  [\"71-37=34\", \"59+37=96\"],\n  [\"1+3=4\", \"33+8=41\"],\n  [\"65-27=38\", \"15+29=44\"],\n  [\"18+2=20\", \"26+37=63\"],\n  [\"59+8=67\", \"82+2=84\"],\n  [\"93-46=47\", \"35-34=1\"],\n  [\"32-4=28\", \"10+71=81\"],\n  [\"2+97=99\", \"86-42=44\"],\n  [\"62-16=46\", \"2+58=60\"],\n  [\"10+46=56\", \"92-68=24\"],\n  [\"87-57=30\", \"1+75=76\"],\n  [\"4+87=91\", \"0+43=43\"],\n  [\"79-5=74\", \"75-39=36\"],\n  [\"97-28=69\", \"82-33=49\"],\n  [\"67-17=50\", \"6-1=5\"],\n  [\"99-30=69\", \"54+7=61\"],\n  [\"82-72=10\", \"46-37=9\"],\n  [\"97-96=1\", \"59+23=82\"],\n  [\"11+60=71\", \"77-12=65\"],\n  [\"66-52=14\", \"95-30=65\"],\n  [\"68-5=63\", \"91-72=19\"],\n  [\"78-19=59\", \"10+86=96\"],\n  [\"36+46=82\", \"95-36=59\"],\n  [\"74-34=40\", \"56-3=53\"],\n  [\"91-45=46\", \"61+10=71\"],\n  [\"93-40=53\", \"74-63=11\"],\n  [\"22+43=65\", \"19-14=5\"],\n  [\"11-9=2\", \"16+58=74\"],\n  [\"91-8=83\", \"20+65=85\"],\n  [\"26+45=71\", \"69-4=65\"],\n  [\"24-7=17\", \"94-1=93\"],\n  [\"54+24=78\", \"39+29=68\"],\n  [\"80-38=42\", \"5+73=78\"],\n  [\"64-6=58\", \"23+7=30\"],\n  [\"22-19=3\", \"18+76=94\"],\n  [\"51+20=71\", \"4+11=15\"],\n  [\"25+26=51\", \"62-21=41\"],\n  [\"84-16=68\", \"49-8=41\"],\n  [\"35+61=96\", \"57-30=27\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\" but found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2023-05-01 Monday\", \"2023-05-02 Tuesday\")\n    ,@(\"84+10=94\", \"50+0=50\")\n    ,@(\"56-36=20\", \"4+60=64\")\n    ,@(\"9+5=14\", \"36-17=19\")\n    ,@(\"38+3=41\", \"77+18=95\")\n    ,@(\"85-0=85\", \"33+0=33\")\n    ,@(\"50-39=11\", \"70+18=88\")\n    ,@(\"0+58=58\", \"85-59=26\")\n    ,@(\"68-11=57\", \"59-29=30\")\n    ,@(\"23+47=70\", \"72-28=44\")\n    ,@(\"44+11=55\", \"10+62=72\")\n    ,@(\"79-0=79\", \"83-9=74\")\n    ,@(\"25+73=98\", \"39+50=89\")\n    ,@(\"46-31=15\", \"67+14=81\")\n    ,@(\"27+46=73\", \"76+22=98\")\n    ,@(\"73-52=21\", \"47-18=29\")\n    ,@(\"94-60=34\", \"59-45=14\")\n    ,@(\"40-0=40\", \"80-14=66\")\n    ,@(\"2+69=71\", \"81-15=66\")\n    ,@(\"15+56=71\", \"58-14=44\")\n    ,@(\"8+56=64\", \"65-16=49\")\n    ,@(\"58-57=1\", \"86-19=67\")\n    ,@(\"9-6=3\", \"40+43=83\")\n    ,@(\"10+43=53\", \"68-24=44\")\n    ,@(\"36+25=61\", \"54+13=67\")\n    ,@(\"85+9=94\", \"28+40=68\")\n    ,@(\"71-65=6\", \"1+88=89\")\n    ,@(\"86-35=51\", \"0+16=16\")\n    ,@(\"94-57=37\", \"44-35=9\")\n    ,@(\"83-26=57\", \"94+0=94\")\n    ,@(\"70-15=55\", \"3+86=89\")\n    ,@(\"70-38=32\", \"49-32=17\")\n    ,@(\"50+14=64\", \"64-51=13\")\n    ,@(\"79-15=64\", \"25+44=69\")\n    ,@(\"97-37=60\", \"71-54=17\")\n    ,@(\"42+50=92\", \"16+71=87\")\n    ,@(\"36+8=44\", \"30+25=55\")\n    ,@(\"39+24=63\", \"36-33=3\")\n    ,@(\"13+34=47\", \"98-55=43\")\n    ,@(\"31+12=43\", \"64-43=21\")\n    ,@(\"44-9=35\", \"74-11=63\")\n    ,@(\"7+53=60\", \"43-12=31\")\n    ,@(\"68-31=37\", \"47+46=93\")\n    ,@(\"72+19=91\", \"86-2=84\")\n    ,@(\"39+31=70\", \"64-63=1\")\n    ,@(\"53-34=19\", \"43+34=77\")\n    ,@(\"34+43=77\", \"55-12=43\")\n    ,@(\"83+3=86\", \"33-23=10\")\n    ,@(\"86-70=16\", \"19-19=0\")\n    ,@(\"36-3=33\", \"33+8=41\")\n    ,@(\"45-32=13\", \"74-3=71\")\n    ,@(\"55-50=5\", \"10+21=31\")\n    ,@(\"11+57=68\", \"59-49=10\")\n    ,@(\"49-44=5\", \"58-51=7\")\n    ,@(\"45-28=17\", \"87-86=1\")\n    ,@(\"82-52=30\", \"54-49=5\")\n    ,@(\"55-46=9\", \"28+47=75\")\n    ,@(\"61-29=32\", \"96-18=78\")\n    ,@(\"11+5=16\", \"43+18=61\")\n    ,@(\"26-24=2\", \"85-19=66\")\n    ,@(\"44-27=17\", \"66-19=47\")\n    ,@(\"12+62=74\", \"81-45=36\")\n    ,@(\"71-37=34\", \"59+37=96\")\n    ,@(\"1+3=4\", \"33+8=41\")\n    ,@(\"65-27=38\", \"15+29=44\")\n    ,@(\"18+2=20\", \"26+37=63\")\n    ,@(\"59+8=67\", \"82+2=84\")\n    ,@(\"93-46=47\", \"35-34=1\")\n    ,@(\"32-4=28\", \"10+71=81\")\n    ,@(\"2+97=99\", \"86-42=44\")\n    ,@(\"62-16=46\", \"2+58=60\")\n    ,@(\"10+46=56\", \"92-68=24\")\n    ,@(\"87-57=30\", \"1+75=76\")\n    ,@(\"4+87=91\", \"0+43=43\")\n    ,@(\"79-5=74\", \"75-39=36\")\n    ,@(\"97-28=69\", \"82-33=49\")\n    ,@(\"67-17=50\", \"6-1=5\")\n    ,@(\"99-30=69\", \"54+7=61\")\n    ,@(\"82-72=10\", \"46-37=9\")\n    ,@(\"97-96=1\", \"59+23=82\")\n    ,@(\"11+60=71\", \"77-12=65\")\n    ,@(\"66-52=14\", \"95-30=65\")\n    ,@(\"68-5=63\", \"91-72=19\")\n    ,@(\"78-19=59\", \"10+86=96\")\n    ,@(\"36+46=82\", \"95-36=59\")\n    ,@(\"74-34=40\", \"56-3=53\")\n    ,@(\"91-45=46\", \"61+10=71\")\n    ,@(\"93-40=53\", \"74-63=11\")\n    ,@(\"22+43=65\", \"19-14=5\")\n    ,@(\"11-9=2\", \"16+58=74\")\n    ,@(\"91-8=83\", \"20+65=85\")\n    ,@(\"26+45=71\", \"69-4=65\")\n    ,@(\"24-7=17\", \"94-1=93\")\n    ,@(\"54+24=78\", \"39+29=68\")\n    ,@(\"80-38=42\", \"5+73=78\")\n    ,@(\"64-6=58\", \"23+7=30\")\n    ,@(\"22-19=3\", \"18+76=94\")\n    ,@(\"51+20=71\", \"4+11=15\")\n    ,@(\"25+26=51\", \"62-21=41\")\n    ,@(\"84-16=68\", \"49-8=41\")\n    ,@(\"35+61=96\", \"57-30=27\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"Find/Replace failed for: $oldText\"\n    }\n}"}
